$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1569.3125
$ws.Range("I40").Value = 1261
$ws.Range("J40").Value = 2083.1667
$ws.Range("K40").Value = 1261
$ws.Range("L40").Value = 2083.1667
$ws.Range("M40").Value = -1086
$ws.Range("N40").Value = -2433.1667
$ws.Range("H52").Value = 2665
$ws.Range("J52").Value = 2898
$ws.Range("L52").Value = 8694
$ws.Range("N52").Value = -9014
$ws.Range("H64").Value = 3892.8
$ws.Range("I64").Value = 3999
$ws.Range("J64").Value = 3866.25
$ws.Range("K64").Value = 3999
$ws.Range("L64").Value = 3866.25
$ws.Range("M64").Value = -3751
$ws.Range("N64").Value = -4362.25
$ws.Range("H67").Value = 3892.8
$ws.Range("I67").Value = 3999
$ws.Range("J67").Value = 3866.25
$ws.Range("K67").Value = 3999
$ws.Range("L67").Value = 3866.25
$ws.Range("M67").Value = -3141
$ws.Range("N67").Value = -5582.25
$ws.Range("H96").Value = 1638.3334
$ws.Range("I96").Value = 1638.3334
$ws.Range("K96").Value = 4915.0002
$ws.Range("M96").Value = -3542.0002
$ws.Range("H137").Value = 14785
$ws.Range("I137").Value = 14082.667
$ws.Range("K137").Value = 42248.001
$ws.Range("M137").Value = -39698.001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 250
$ws.Range("I5").Value = 250
$ws.Range("K5").Value = 250
$ws.Range("M5").Value = -138
$ws.Range("H41").Value = 2080.5
$ws.Range("I41").Value = 2107.3333
$ws.Range("J41").Value = 2000
$ws.Range("K41").Value = 2107.3333
$ws.Range("L41").Value = 2000
$ws.Range("M41").Value = -1693.3333
$ws.Range("N41").Value = -2828
$ws.Range("H45").Value = 2526.25
$ws.Range("I45").Value = 2114.8572
$ws.Range("J45").Value = 3102.2
$ws.Range("K45").Value = 2114.8572
$ws.Range("L45").Value = 3102.2
$ws.Range("M45").Value = -1737.8572
$ws.Range("N45").Value = -3856.2
$ws.Range("H74").Value = 3355.1667
$ws.Range("I74").Value = 3521.6667
$ws.Range("K74").Value = 3521.6667
$ws.Range("M74").Value = -2647.6667
$ws.Range("H77").Value = 3355.1667
$ws.Range("I77").Value = 3521.6667
$ws.Range("K77").Value = 17608.3335
$ws.Range("M77").Value = -13240.3335
$ws.Range("H122").Value = 4352.4165
$ws.Range("I122").Value = 3746.5
$ws.Range("K122").Value = 11239.5
$ws.Range("M122").Value = -8789.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 250
$ws.Range("I4").Value = 250
$ws.Range("K4").Value = 250
$ws.Range("M4").Value = -135
$ws.Range("H22").Value = 287
$ws.Range("I22").Value = 287
$ws.Range("K22").Value = 287
$ws.Range("M22").Value = -114
$ws.Range("H99").Value = 1951.8667
$ws.Range("I99").Value = 1614.8334
$ws.Range("J99").Value = 3300
$ws.Range("K99").Value = 1614.8334
$ws.Range("L99").Value = 3300
$ws.Range("M99").Value = -116.8334
$ws.Range("N99").Value = -6296
$ws.Range("H107").Value = 1307.5217
$ws.Range("I107").Value = 801.5789
$ws.Range("K107").Value = 801.5789
$ws.Range("M107").Value = 1118.4211

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 4733.3335
$ws.Range("H16").Value = 1138.5294
$ws.Range("I16").Value = 1273.5555
$ws.Range("J16").Value = 986.625
$ws.Range("K16").Value = 1273.5555
$ws.Range("L16").Value = 986.625
$ws.Range("M16").Value = -986.5554999999999
$ws.Range("N16").Value = -1560.625
$ws.Range("H69").Value = 17937.375
$ws.Range("I69").Value = 16214.143
$ws.Range("K69").Value = 16214.143
$ws.Range("M69").Value = -15465.143
$ws.Range("H72").Value = 17937.375
$ws.Range("I72").Value = 16214.143
$ws.Range("K72").Value = 48642.429
$ws.Range("M72").Value = -44898.429
$ws.Range("H99").Value = 13434.759
$ws.Range("I99").Value = 8760.071
$ws.Range("K99").Value = 8760.071
$ws.Range("M99").Value = -7262.071
$ws.Range("H107").Value = 1066.7838
$ws.Range("I107").Value = 818.9091
$ws.Range("J107").Value = 1171.6538
$ws.Range("K107").Value = 818.9091
$ws.Range("L107").Value = 1171.6538
$ws.Range("M107").Value = 1101.0909
$ws.Range("N107").Value = -5011.6538
$ws.Range("H113").Value = 1138.5294
$ws.Range("I113").Value = 1273.5555
$ws.Range("J113").Value = 986.625
$ws.Range("K113").Value = 1273.5555
$ws.Range("L113").Value = 986.625
$ws.Range("M113").Value = 896.4445000000001
$ws.Range("N113").Value = -5326.625
$ws.Range("H122").Value = 1837.75
$ws.Range("I122").Value = 1750.6
$ws.Range("K122").Value = 5251.799999999999
$ws.Range("M122").Value = -2801.799999999999
$ws.Range("H126").Value = 13434.759
$ws.Range("I126").Value = 8760.071
$ws.Range("K126").Value = 26280.213
$ws.Range("M126").Value = -23810.213
$ws.Range("H132").Value = 4211.8335
$ws.Range("I132").Value = 3268.1667
$ws.Range("K132").Value = 9804.500100000001
$ws.Range("M132").Value = -7274.500100000001

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H122").Value = 200
$ws.Range("J122").Value = 200
$ws.Range("L122").Value = 1800
$ws.Range("N122").Value = -6700
$ws.Range("H128").Value = 953289.8
$ws.Range("I128").Value = 953289.8
$ws.Range("K128").Value = 2859869.4
$ws.Range("M128").Value = -2854889.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2446.4285
$ws.Range("I80").Value = 2333.3333
$ws.Range("J80").Value = 2650
$ws.Range("K80").Value = 2333.3333
$ws.Range("L80").Value = 2650
$ws.Range("M80").Value = -1335.3333
$ws.Range("N80").Value = -4646
$ws.Range("H83").Value = 2446.4285
$ws.Range("I83").Value = 2333.3333
$ws.Range("J83").Value = 2650
$ws.Range("K83").Value = 11666.6665
$ws.Range("L83").Value = 13250
$ws.Range("M83").Value = -6674.666499999999
$ws.Range("N83").Value = -23234
$ws.Range("H102").Value = 2063.0952
$ws.Range("I102").Value = 566.0714
$ws.Range("K102").Value = 566.0714
$ws.Range("M102").Value = 1055.9286

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 18000000
$ws.Range("J2").Value = 15999999
$ws.Range("L2").Value = 15999999
$ws.Range("N2").Value = -16000223
$ws.Range("H16").Value = 10014.167
$ws.Range("I16").Value = 8699.25
$ws.Range("K16").Value = 8699.25
$ws.Range("M16").Value = -8529.25
$ws.Range("H46").Value = 1045.8
$ws.Range("I46").Value = 1045.8
$ws.Range("J46").Value = 0
$ws.Range("K46").Value = 1045.8
$ws.Range("L46").Value = 0
$ws.Range("M46").Value = -857.8
$ws.Range("N46").ClearContents()
$ws.Range("H55").Value = 1074.2307
$ws.Range("I55").Value = 1010.7778
$ws.Range("J55").Value = 1217
$ws.Range("K55").Value = 1010.7778
$ws.Range("L55").Value = 1217
$ws.Range("M55").Value = -837.7778
$ws.Range("N55").Value = -1563
$ws.Range("H68").Value = 3962.0625
$ws.Range("J68").Value = 5314.6665
$ws.Range("L68").Value = 5314.6665
$ws.Range("N68").Value = -6812.6665
$ws.Range("H71").Value = 3962.0625
$ws.Range("J71").Value = 5314.6665
$ws.Range("L71").Value = 26573.3325
$ws.Range("N71").Value = -34061.3325
$ws.Range("H82").Value = 168216
$ws.Range("I82").Value = 2136.75
$ws.Range("K82").Value = 2136.75
$ws.Range("M82").Value = -1775.75
$ws.Range("H85").Value = 168216
$ws.Range("I85").Value = 2136.75
$ws.Range("K85").Value = 2136.75
$ws.Range("M85").Value = -888.75
$ws.Range("H93").Value = 6399.6
$ws.Range("I93").Value = 7249.5
$ws.Range("J93").Value = 3000
$ws.Range("K93").Value = 7249.5
$ws.Range("L93").Value = 3000
$ws.Range("M93").Value = -6001.5
$ws.Range("N93").Value = -5496
$ws.Range("H100").Value = 10400.2
$ws.Range("I100").Value = 6750.75
$ws.Range("K100").Value = 6750.75
$ws.Range("M100").Value = -6209.75
$ws.Range("H122").Value = 3283.2856
$ws.Range("I122").Value = 2897.25
$ws.Range("K122").Value = 8691.75
$ws.Range("M122").Value = -6241.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("I2").Value = 667000
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 667000
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -666888
$ws.Range("N2").ClearContents()
$ws.Range("H41").Value = 147739.5
$ws.Range("J41").Value = 19991
$ws.Range("L41").Value = 19991
$ws.Range("N41").Value = -20771
$ws.Range("H119").Value = 15475
$ws.Range("I119").Value = 950
$ws.Range("J119").Value = 30000
$ws.Range("K119").Value = 950
$ws.Range("L119").Value = 30000
$ws.Range("M119").Value = 3888
$ws.Range("N119").Value = -39676
$ws.Range("H122").Value = 1508.5385
$ws.Range("I122").Value = 1534.25
$ws.Range("K122").Value = 4602.75
$ws.Range("M122").Value = -2152.75
$ws.Range("H132").Value = 3663.2856
$ws.Range("I132").Value = 3913
$ws.Range("J132").Value = 3476
$ws.Range("K132").Value = 11739
$ws.Range("L132").Value = 10428
$ws.Range("M132").Value = -9209
$ws.Range("N132").Value = -15488
$ws.Range("H136").Value = 11285.857
$ws.Range("I136").Value = 12832.667
$ws.Range("J136").Value = 2005
$ws.Range("K136").Value = 38498.001
$ws.Range("L136").Value = 6015
$ws.Range("M136").Value = -35948.001
$ws.Range("N136").Value = -11115
